$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "fhgutyr"
$ws.Range("B18").Value = "hgf"
$ws.Range("C18").Value = "gf"
$ws.Range("D18").Value = ""
